$d = $word.ActiveDocument
$full = $d.Content.WordOpenXML()

$startMarker = '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>'
$startIdx = $full.IndexOf($startMarker)
$contentStart = $startIdx + $startMarker.Length
$endMarker = '</pkg:xmlData></pkg:part>'
$endIdx = $full.IndexOf($endMarker, $contentStart)

$docXml = $full.Substring($contentStart, $endIdx - $contentStart)

$d.Content.InsertXML($docXml)
